$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly record at row 9, pushing the existing rows 9:15 down to 10:16.
# Row 9 currently holds the most-recent-so-far record; a newer record needs to be
# inserted above it, so first shift the existing data block down by one row.
$src = $ws.Range("A9:R15")
$dst = $ws.Range("A10:R16")
$src.Copy($dst)

# Now overwrite row 9 with the new weekly record's values.
$ws.Cells.Item(9, 4).Value2 = 44740   # D9  Fecha
$ws.Cells.Item(9, 10).Value2 = 120    # J9  Volumen
$ws.Cells.Item(9, 11).Value2 = 6000   # K9  Precio minimo
$ws.Cells.Item(9, 12).Value2 = 7000   # L9  Precio maximo
$ws.Cells.Item(9, 13).Value2 = 6500   # M9  Precio promedio ponderado
$ws.Cells.Item(9, 16).Value2 = 108    # P9  Precio $/Kg
